$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new text value, and whether the value "looks numeric"
# (bare price figures like 0.5120/1.001/15.11) and therefore needs to be forced back
# to Text so Excel does not silently reinterpret it as a Number and drop formatting
# (e.g. trailing zeros, like "0.5120" -> 0.512).
$updates = @(
    @{ Cell = 'D2'; Value = '25.895.53'; ForceText = $false }
    @{ Cell = 'E2'; Value = '  -0.71%  '; ForceText = $false }
    @{ Cell = 'D3'; Value = '1.740.81'; ForceText = $false }
    @{ Cell = 'E3'; Value = '  -0.50%  '; ForceText = $false }
    @{ Cell = 'E4'; Value = '  +0.20%  '; ForceText = $false }
    @{ Cell = 'D5'; Value = '248.43'; ForceText = $true }
    @{ Cell = 'E5'; Value = '  +5.85%  '; ForceText = $false }
    @{ Cell = 'E6'; Value = '  +0.15%  '; ForceText = $false }
    @{ Cell = 'D7'; Value = '0.5120'; ForceText = $true }
    @{ Cell = 'E7'; Value = '  -3.20%  '; ForceText = $false }
    @{ Cell = 'B8'; Value = 'Cardano'; ForceText = $false }
    @{ Cell = 'C8'; Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'; ForceText = $false }
    @{ Cell = 'D8'; Value = '0.2743'; ForceText = $true }
    @{ Cell = 'E8'; Value = '  -1.92%  '; ForceText = $false }
    @{ Cell = 'B9'; Value = 'Dogecoin'; ForceText = $false }
    @{ Cell = 'C9'; Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'; ForceText = $false }
    @{ Cell = 'D9'; Value = '0.06180'; ForceText = $true }
    @{ Cell = 'E9'; Value = '  -0.07%  '; ForceText = $false }
    @{ Cell = 'B10'; Value = 'WrappedEther'; ForceText = $false }
    @{ Cell = 'C10'; Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; ForceText = $false }
    @{ Cell = 'D10'; Value = '1.741.21'; ForceText = $false }
    @{ Cell = 'E10'; Value = '  -0.21%  '; ForceText = $false }
    @{ Cell = 'B11'; Value = 'TRON'; ForceText = $false }
    @{ Cell = 'C11'; Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; ForceText = $false }
    @{ Cell = 'D11'; Value = '0.07237'; ForceText = $true }
    @{ Cell = 'E11'; Value = '  +0.68%  '; ForceText = $false }
    @{ Cell = 'B12'; Value = 'Solana'; ForceText = $false }
    @{ Cell = 'C12'; Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; ForceText = $false }
    @{ Cell = 'D12'; Value = '15.11'; ForceText = $true }
    @{ Cell = 'E12'; Value = '  -1.90%  '; ForceText = $false }
    @{ Cell = 'B13'; Value = 'Polygon'; ForceText = $false }
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'; ForceText = $false }
    @{ Cell = 'D13'; Value = '0.6482'; ForceText = $true }
    @{ Cell = 'E13'; Value = '  +0.51%  '; ForceText = $false }
    @{ Cell = 'B14'; Value = 'Polkadot'; ForceText = $false }
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'; ForceText = $false }
    @{ Cell = 'D14'; Value = '4.626'; ForceText = $true }
    @{ Cell = 'E14'; Value = '  +0.20%  '; ForceText = $false }
    @{ Cell = 'B15'; Value = 'Litecoin'; ForceText = $false }
    @{ Cell = 'C15'; Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; ForceText = $false }
    @{ Cell = 'D15'; Value = '77.64'; ForceText = $true }
    @{ Cell = 'E15'; Value = '  -0.95%  '; ForceText = $false }
    @{ Cell = 'B16'; Value = 'Dai'; ForceText = $false }
    @{ Cell = 'C16'; Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'; ForceText = $false }
    @{ Cell = 'D16'; Value = '1.001'; ForceText = $true }
    @{ Cell = 'E16'; Value = '  +0.14%  '; ForceText = $false }
    @{ Cell = 'B17'; Value = 'BinanceUSD'; ForceText = $false }
    @{ Cell = 'C17'; Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; ForceText = $false }
    @{ Cell = 'E17'; Value = '  +0.19%  '; ForceText = $false }
    @{ Cell = 'B18'; Value = 'WrappedBTC'; ForceText = $false }
    @{ Cell = 'C18'; Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; ForceText = $false }
    @{ Cell = 'D18'; Value = '25.925.03'; ForceText = $false }
    @{ Cell = 'E18'; Value = '  -0.26%  '; ForceText = $false }
    @{ Cell = 'B19'; Value = 'Avalanche'; ForceText = $false }
    @{ Cell = 'C19'; Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; ForceText = $false }
    @{ Cell = 'D19'; Value = '11.82'; ForceText = $true }
    @{ Cell = 'E19'; Value = '  +1.25%  '; ForceText = $false }
    @{ Cell = 'B20'; Value = 'ShibaInu'; ForceText = $false }
    @{ Cell = 'C20'; Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; ForceText = $false }
    @{ Cell = 'D20'; Value = '0.000006807'; ForceText = $true }
    @{ Cell = 'E20'; Value = '  +1.26%  '; ForceText = $false }
    @{ Cell = 'B21'; Value = 'WrappedliquidstakedEther2.0'; ForceText = $false }
    @{ Cell = 'C21'; Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; ForceText = $false }
    @{ Cell = 'D21'; Value = '1.966.59'; ForceText = $false }
    @{ Cell = 'E21'; Value = '  -0.16%  '; ForceText = $false }
    @{ Cell = 'B22'; Value = 'Uniswap'; ForceText = $false }
    @{ Cell = 'C22'; Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; ForceText = $false }
    @{ Cell = 'D22'; Value = '4.277'; ForceText = $true }
    @{ Cell = 'E22'; Value = '  -0.70%  '; ForceText = $false }
    @{ Cell = 'B23'; Value = 'Cosmos'; ForceText = $false }
    @{ Cell = 'C23'; Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; ForceText = $false }
    @{ Cell = 'D23'; Value = '8.648'; ForceText = $true }
    @{ Cell = 'E23'; Value = '  -1.33%  '; ForceText = $false }
    @{ Cell = 'B24'; Value = 'Chainlink'; ForceText = $false }
    @{ Cell = 'C24'; Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; ForceText = $false }
    @{ Cell = 'D24'; Value = '5.390'; ForceText = $true }
    @{ Cell = 'E24'; Value = '  +3.12%  '; ForceText = $false }
    @{ Cell = 'B25'; Value = 'Monero'; ForceText = $false }
    @{ Cell = 'C25'; Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; ForceText = $false }
    @{ Cell = 'D25'; Value = '135.79'; ForceText = $true }
    @{ Cell = 'E25'; Value = '  -2.01%  '; ForceText = $false }
    @{ Cell = 'B26'; Value = 'Toncoin'; ForceText = $false }
    @{ Cell = 'C26'; Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; ForceText = $false }
    @{ Cell = 'D26'; Value = '1.503'; ForceText = $true }
    @{ Cell = 'E26'; Value = '  -0.16%  '; ForceText = $false }
    @{ Cell = 'B27'; Value = 'EthereumClassic'; ForceText = $false }
    @{ Cell = 'C27'; Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; ForceText = $false }
    @{ Cell = 'D27'; Value = '15.21'; ForceText = $true }
    @{ Cell = 'E27'; Value = '  -0.49%  '; ForceText = $false }
    @{ Cell = 'B28'; Value = 'LidoDAOToken'; ForceText = $false }
    @{ Cell = 'C28'; Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; ForceText = $false }
    @{ Cell = 'D28'; Value = '1.771'; ForceText = $true }
    @{ Cell = 'E28'; Value = '  -1.80%  '; ForceText = $false }
    @{ Cell = 'B29'; Value = 'BitcoinCash'; ForceText = $false }
    @{ Cell = 'C29'; Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; ForceText = $false }
    @{ Cell = 'D29'; Value = '105.80'; ForceText = $true }
    @{ Cell = 'E29'; Value = '  +1.18%  '; ForceText = $false }
    @{ Cell = 'B30'; Value = 'InternetComputer(DFINITY)'; ForceText = $false }
    @{ Cell = 'C30'; Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; ForceText = $false }
    @{ Cell = 'D30'; Value = '3.913'; ForceText = $true }
    @{ Cell = 'E30'; Value = '  +3.07%  '; ForceText = $false }
    @{ Cell = 'B31'; Value = 'Stellar'; ForceText = $false }
    @{ Cell = 'C31'; Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; ForceText = $false }
    @{ Cell = 'D31'; Value = '0.08226'; ForceText = $true }
    @{ Cell = 'E31'; Value = '  -0.72%  '; ForceText = $false }
    @{ Cell = 'B32'; Value = 'Filecoin'; ForceText = $false }
    @{ Cell = 'C32'; Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; ForceText = $false }
    @{ Cell = 'D32'; Value = '3.644'; ForceText = $true }
    @{ Cell = 'E32'; Value = '  -0.71%  '; ForceText = $false }
    @{ Cell = 'B33'; Value = 'Hedera'; ForceText = $false }
    @{ Cell = 'C33'; Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'; ForceText = $false }
    @{ Cell = 'D33'; Value = '0.04695'; ForceText = $true }
    @{ Cell = 'E33'; Value = '  +2.98%  '; ForceText = $false }
    @{ Cell = 'B34'; Value = 'HuobiToken'; ForceText = $false }
    @{ Cell = 'C34'; Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; ForceText = $false }
    @{ Cell = 'D34'; Value = '2.653'; ForceText = $true }
    @{ Cell = 'E34'; Value = '  +0.36%  '; ForceText = $false }
    @{ Cell = 'B35'; Value = 'ARBITRUM'; ForceText = $false }
    @{ Cell = 'C35'; Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; ForceText = $false }
    @{ Cell = 'D35'; Value = '0.9983'; ForceText = $true }
    @{ Cell = 'E35'; Value = '  -0.50%  '; ForceText = $false }
    @{ Cell = 'B36'; Value = 'ImmutableX'; ForceText = $false }
    @{ Cell = 'C36'; Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; ForceText = $false }
    @{ Cell = 'D36'; Value = '0.6260'; ForceText = $true }
    @{ Cell = 'E36'; Value = '  -1.31%  '; ForceText = $false }
    @{ Cell = 'B37'; Value = 'MXToken'; ForceText = $false }
    @{ Cell = 'C37'; Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; ForceText = $false }
    @{ Cell = 'D37'; Value = '2.727'; ForceText = $true }
    @{ Cell = 'E37'; Value = '  +0.68%  '; ForceText = $false }
    @{ Cell = 'B38'; Value = 'VeChain'; ForceText = $false }
    @{ Cell = 'C38'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; ForceText = $false }
    @{ Cell = 'D38'; Value = '0.01603'; ForceText = $true }
    @{ Cell = 'E38'; Value = '  +0.63%  '; ForceText = $false }
    @{ Cell = 'B39'; Value = 'RenderToken'; ForceText = $false }
    @{ Cell = 'C39'; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; ForceText = $false }
    @{ Cell = 'D39'; Value = '1.914'; ForceText = $true }
    @{ Cell = 'E39'; Value = '  -1.17%  '; ForceText = $false }
    @{ Cell = 'B40'; Value = 'PaxDollar'; ForceText = $false }
    @{ Cell = 'C40'; Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'; ForceText = $false }
    @{ Cell = 'D40'; Value = '1.000'; ForceText = $true }
    @{ Cell = 'E40'; Value = '  +0.15%  '; ForceText = $false }
    @{ Cell = 'B41'; Value = 'Quant'; ForceText = $false }
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; ForceText = $false }
    @{ Cell = 'D41'; Value = '100.07'; ForceText = $true }
    @{ Cell = 'E41'; Value = '  +0.94%  '; ForceText = $false }
    @{ Cell = 'B42'; Value = 'TrustWalletToken'; ForceText = $false }
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; ForceText = $false }
    @{ Cell = 'D42'; Value = '0.7573'; ForceText = $true }
    @{ Cell = 'E42'; Value = '  +1.83%  '; ForceText = $false }
    @{ Cell = 'B43'; Value = 'TheSandbox'; ForceText = $false }
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; ForceText = $false }
    @{ Cell = 'D43'; Value = '0.3840'; ForceText = $true }
    @{ Cell = 'E43'; Value = '  -2.12%  '; ForceText = $false }
    @{ Cell = 'B44'; Value = 'FraxShare'; ForceText = $false }
    @{ Cell = 'C44'; Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; ForceText = $false }
    @{ Cell = 'D44'; Value = '4.994'; ForceText = $true }
    @{ Cell = 'E44'; Value = '  -0.61%  '; ForceText = $false }
    @{ Cell = 'B45'; Value = 'Algorand'; ForceText = $false }
    @{ Cell = 'C45'; Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; ForceText = $false }
    @{ Cell = 'D45'; Value = '0.1130'; ForceText = $true }
    @{ Cell = 'E45'; Value = '  -1.63%  '; ForceText = $false }
    @{ Cell = 'B46'; Value = 'Aptos'; ForceText = $false }
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; ForceText = $false }
    @{ Cell = 'D46'; Value = '6.288'; ForceText = $true }
    @{ Cell = 'E46'; Value = '  -0.69%  '; ForceText = $false }
    @{ Cell = 'B47'; Value = 'Aave'; ForceText = $false }
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'; ForceText = $false }
    @{ Cell = 'D47'; Value = '55.28'; ForceText = $true }
    @{ Cell = 'E47'; Value = '  +2.15%  '; ForceText = $false }
    @{ Cell = 'B48'; Value = 'Cronos'; ForceText = $false }
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; ForceText = $false }
    @{ Cell = 'D48'; Value = '0.05232'; ForceText = $true }
    @{ Cell = 'E48'; Value = '  -2.08%  '; ForceText = $false }
    @{ Cell = 'B49'; Value = 'Elrond'; ForceText = $false }
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'; ForceText = $false }
    @{ Cell = 'D49'; Value = '30.67'; ForceText = $true }
    @{ Cell = 'E49'; Value = '  -0.39%  '; ForceText = $false }
    @{ Cell = 'B50'; Value = 'EnergySwap'; ForceText = $false }
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; ForceText = $false }
    @{ Cell = 'D50'; Value = '7.491'; ForceText = $true }
    @{ Cell = 'E50'; Value = '  -2.18%  '; ForceText = $false }
    @{ Cell = 'D51'; Value = '0.3409'; ForceText = $true }
    @{ Cell = 'E51'; Value = '  -1.44%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Apply a Text number format so the assigned string survives verbatim,
        # then restore the default style so no visible formatting change remains.
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}
